# Applies the "Natmi following Dr Hou advice" data refresh to Sheet1 (A2:T16).
# Source data recomputed with updated cluster counts (per commit); this script
# rewrites every data cell, which also naturally extends the sheet from 12 to 15
# data rows (dimension grows from A1:T13 to A1:T16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Cthrc1", "Fzd5", "ECs", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 2, 0.6666666666666666, 2.676920666666666, 8.030762, 0.1137411923116975, 0.1167596329733683, 1.657290507802222, 14.91561457022, 0.00090192645930242, 0.0009296974251082671)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $row2[$i] }

$row3 = @("ECs", "Cthrc1", "Fzd5", "FAPs", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 3, 1, 10.75807, 32.27421, 0.4571057050773153, 0.4692362834442626, 6.660356997233333, 59.9432129751, 0.003624682682923832, 0.003736289275463958)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $row3[$i] }

$row4 = @("ECs", "Cthrc1", "Fzd5", "M1", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 3, 1, 3.564735333333333, 10.694206, 0.1514640505181089, 0.1554835727296604, 2.206939527317778, 19.86245574586, 0.001201055062101292, 0.001238036413204299)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $row4[$i] }

$row5 = @("ECs", "Cthrc1", "Fzd5", "M2", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 3, 1, 4.710187333333333, 14.130562, 0.2001338067190093, 0.2054449170361947, 2.916092678691111, 26.24483410822, 0.001586988601158062, 0.001635853124116083)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $row5[$i] }

$row6 = @("FAPs", "Cthrc1", "Fzd5", "sCs", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 2, 1, 1.8252775, 3.650555, 0.07755524537386906, 0.05307559381651387, 1.130035384508333, 6.780212307049999, 0.0006149850062121278, 0.000422613891896698)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i + 1).Value = $row6[$i] }

$row7 = @("FAPs", "Cthrc1", "Fzd5", "ECs", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 2, 0.6666666666666666, 2.676920666666666, 8.030762, 0.1137411923116975, 0.1167596329733683, 204.7553069954535, 1842.797762959082, 0.1114314166359952, 0.1148624702758771)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i + 1).Value = $row7[$i] }

$row8 = @("FAPs", "Cthrc1", "Fzd5", "FAPs", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 3, 1, 10.75807, 32.27421, 0.4571057050773153, 0.4692362834442626, 822.8753107844232, 7405.877797059808, 0.4478231257641059, 0.4616119226049052)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i + 1).Value = $row8[$i] }

$row9 = @("FAPs", "Cthrc1", "Fzd5", "M1", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 3, 1, 3.564735333333333, 10.694206, 0.1514640505181089, 0.1554835727296604, 272.6634698678184, 2453.971228810366, 0.1483882257221867, 0.1529572061529287)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i + 1).Value = $row9[$i] }

$row10 = @("sCs", "Cthrc1", "Fzd5", "M2", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 3, 1, 4.710187333333333, 14.130562, 0.2001338067190093, 0.2054449170361947, 360.2780857318758, 3242.502771586881, 0.1960696309419656, 0.2021067562089921)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i + 1).Value = $row10[$i] }

$row11 = @("sCs", "Cthrc1", "Fzd5", "sCs", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 2, 1, 1.8252775, 3.650555, 0.07755524537386906, 0.05307559381651387, 139.6138703392258, 837.6832220353549, 0.07598030831151804, 0.05221319784821843)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, $i + 1).Value = $row11[$i] }

$row12 = @("sCs", "Cthrc1", "Fzd5", "ECs", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 2, 0.6666666666666666, 2.676920666666666, 8.030762, 0.1137411923116975, 0.1167596329733683, 2.586923932313666, 15.521543593882, 0.001407849216399905, 0.0009674652723830053)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, $i + 1).Value = $row12[$i] }

$row13 = @("sCs", "Cthrc1", "Fzd5", "FAPs", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 3, 1, 10.75807, 32.27421, 0.4571057050773153, 0.4692362834442626, 10.396389065635, 62.37833439380999, 0.00565789663028564, 0.003888071563893477)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, $i + 1).Value = $row13[$i] }

$row14 = @("sCs", "Cthrc1", "Fzd5", "M1", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 3, 1, 3.564735333333333, 10.694206, 0.1514640505181089, 0.1554835727296604, 3.444890713794333, 20.669344282766, 0.00187476973382092, 0.001288330163527442)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, $i + 1).Value = $row14[$i] }

$row15 = @("sCs", "Cthrc1", "Fzd5", "M2", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 3, 1, 4.710187333333333, 14.130562, 0.2001338067190093, 0.2054449170361947, 4.551833190280333, 27.310999141682, 0.002477187175885709, 0.001702307703086574)
for ($i = 0; $i -lt $row15.Length; $i++) { $ws.Cells.Item(15, $i + 1).Value = $row15[$i] }

$row16 = @("sCs", "Cthrc1", "Fzd5", "sCs", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 2, 1, 1.8252775, 3.650555, 0.07755524537386906, 0.05307559381651387, 1.76391258308875, 7.055650332354999, 0.0009599520561388984, 0.0004397820763987455)
for ($i = 0; $i -lt $row16.Length; $i++) { $ws.Cells.Item(16, $i + 1).Value = $row16[$i] }
